# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handoff has been handed back and is in sync with en-US:
#   - Overview sheet: status cells move from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - zh-cn / de-de sheets: add "Latest Target File" (F) and
#     "Latest Handback File" (G) columns (with hyperlinks) for each row,
#     and update the "Latest Handback DateTime" (H) column with the new
#     handback timestamps.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            return $h.Address
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Overview sheet: update status text for both languages / both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper that fills in the new "Latest Target File" / "Latest Handback
# File" columns (F/G) for a language detail sheet, wiring up the
# hyperlinks to the same targets as the existing Source File Name (A)
# and Latest Handoff File (D) hyperlinks, and refreshes the "Latest
# Handback DateTime" (H) column.
# ---------------------------------------------------------------------
function Update-LanguageSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Both rows re-use row 2's "Source File Name" (md) and "Latest
    # Handoff File" (xlf) values/targets for the new "Latest Target
    # File" (F) and "Latest Handback File" (G) columns, matching the
    # source data exactly.
    $mdName = $ws.Range("A2").Value2
    $xlfName = $ws.Range("D2").Value2

    $mdUrl = Get-HyperlinkAddress $ws '$A$2'
    $xlfUrl = Get-HyperlinkAddress $ws '$D$2'

    # Row 2
    $ws.Range("F2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Range("F2").Style = $ws.Range("A2").Style

    $ws.Range("G2").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", $xlfName) | Out-Null
    $ws.Range("G2").Style = $ws.Range("D2").Style

    # Row 3
    $ws.Range("F3").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Range("F3").Style = $ws.Range("A3").Style

    $ws.Range("G3").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, "", "", $xlfName) | Out-Null
    $ws.Range("G3").Style = $ws.Range("D3").Style

    # Latest Handback DateTime
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime
}

Update-LanguageSheet "zh-cn" "2016-03-12 06:43:28"
Update-LanguageSheet "de-de" "2016-03-12 06:43:33"

Write-Host "Handback report generated."
